# Checkout UI: add payment-card columns (CARD NUMBER / CVC / Expiration)
# to the users sheet, used for signup/login while checkout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("P1").Value = "CARD NUMBER"
$ws.Range("Q1").Value = "CVC"
$ws.Range("R1").Value = "Expiration"

# --- Column P: card numbers ------------------------------------------
# (written column-first so the shared-strings table is populated in the
# same P-then-Q-then-R order produced by the original edit)
$ws.Range("P2").Value = "4242 4242 4242 4242"
$ws.Range("P3").Value = "4243 4242 4242 4242"
$ws.Range("P4").Value = "4244 4242 4242 4242"
$ws.Range("P5").Value = "4245 4242 4242 4242"
$ws.Range("P6").Value = "4246 4242 4242 4242"
$ws.Range("P7").Value = "4247 4242 4242 4242"

# --- Column Q: CVC codes ---------------------------------------------
$ws.Range("Q2").Value = "111"
$ws.Range("Q3").Value = "112"
$ws.Range("Q4").Value = "113"
$ws.Range("Q5").Value = "114"
$ws.Range("Q6").Value = "115"
$ws.Range("Q7").Value = "116"

# --- Column R: expiration dates ---------------------------------------
$ws.Range("R2").Value = "12/2020"
$ws.Range("R3").Value = "12/2021"
$ws.Range("R4").Value = "12/2022"
$ws.Range("R5").Value = "12/2023"
$ws.Range("R6").Value = "12/2024"
$ws.Range("R7").Value = "12/2025"

# Reflect the new selection left behind on the sheet after this edit.
$ws.Range("R2:R7").Select()
